$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Insert a new row above the current row 18 (ROCAFUERTE LOPEZ EVELYN ESTEFANIA),
# shifting the existing rows 18-21 down to 19-22 (totals row moves 21 -> 22).
$ws.Rows(18).Insert()

# Populate the newly inserted row 18 with the new sale record. Row-insert
# already carries the neighbouring rows' cell formatting onto the new row,
# so there's no separate style/number-format step needed here.
$ws.Range("A18").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws.Range("B18").Value = "PAUTA ASTUDILLO JULIO HERNAN"
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 326.73
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 1000

# Update the totals row (now row 22) to reflect the new record's contribution.
$ws.Range("D22").Value = 17549.57
$ws.Range("G22").Value = 1000
